# Updated symbol list on Sun Jan 15 06:42:00 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price/volume(1h) snapshot on Sheet1 with the latest
# pull from coinranking.com. Prices and 1h volume deltas are stored as
# plain text (matching the sheet's existing inline-string layout), so each
# assignment is quote-prefixed to keep Excel from re-typing them as
# numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.83"
$ws.Range("E2").Value = "'-3.61%"
$ws.Range("D3").Value = "'31.80"
$ws.Range("E3").Value = "'-0.95%"
$ws.Range("D4").Value = "'5.111"
$ws.Range("E4").Value = "'-4.27%"
$ws.Range("D5").Value = "'0.07521"
$ws.Range("E5").Value = "'-0.10%"
$ws.Range("D6").Value = "'7.755"
$ws.Range("E6").Value = "'-0.66%"
$ws.Range("D7").Value = "'1.717"
$ws.Range("E7").Value = "'9.46%"
$ws.Range("E8").Value = "'3.35%"
$ws.Range("D9").Value = "'0.9292"
$ws.Range("E9").Value = "'2.49%"
$ws.Range("D10").Value = "'0.1698"
$ws.Range("E10").Value = "'0.77%"
$ws.Range("D11").Value = "'0.07499"
$ws.Range("E11").Value = "'-2.27%"
$ws.Range("D12").Value = "'0.07983"
$ws.Range("E12").Value = "'-1.10%"
$ws.Range("D13").Value = "'0.03052"
$ws.Range("E13").Value = "'0.85%"
$ws.Range("D14").Value = "'0.09896"
$ws.Range("E14").Value = "'0.29%"
$ws.Range("D15").Value = "'0.001490"
$ws.Range("E15").Value = "'-2.20%"
$ws.Range("D16").Value = "'0.006495"
$ws.Range("E16").Value = "'0.09%"
$ws.Range("D17").Value = "'3.463"
$ws.Range("E17").Value = "'-1.13%"
$ws.Range("D18").Value = "'2.225"
$ws.Range("E18").Value = "'-0.71%"
$ws.Range("D19").Value = "'0.3280"
$ws.Range("E19").Value = "'0.40%"
$ws.Range("D20").Value = "'0.1325"
$ws.Range("D21").Value = "'4.565"
$ws.Range("E21").Value = "'9.33%"
$ws.Range("D22").Value = "'0.04652"
$ws.Range("E22").Value = "'2.27%"
$ws.Range("D23").Value = "'0.1556"
$ws.Range("E23").Value = "'-4.34%"
$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'0.36%"
$ws.Range("D25").Value = "'0.004424"
$ws.Range("E25").Value = "'-1.43%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'0.14%"
$ws.Range("E27").Value = "'6.84%"
$ws.Range("D39").Value = "'0.01682"
$ws.Range("E39").Value = "'-0.31%"
$ws.Range("D40").Value = "'0.04535"
$ws.Range("E40").Value = "'-0.25%"
$ws.Range("D41").Value = "'0.007083"
$ws.Range("E41").Value = "'-0.99%"
$ws.Range("D42").Value = "'0.1328"
$ws.Range("E42").Value = "'-2.55%"
$ws.Range("D43").Value = "'0.002060"
$ws.Range("E43").Value = "'-8.71%"
$ws.Range("D44").Value = "'0.01169"
$ws.Range("E44").Value = "'-15.78%"
$ws.Range("D45").Value = "'0.00005995"
$ws.Range("E45").Value = "'-1.76%"
$ws.Range("D46").Value = "'1.930"
$ws.Range("E46").Value = "'1.96%"
$ws.Range("E47").Value = "'-0.14%"
